$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray leftover value in C2 (bug: extra forecast value that shouldn't be there)
$ws.Range("C2").ClearContents()

# Recalculated (corrected) values after fixing the naive component forecaster bug.
# Column C (y_0_forecast)
$ws.Range("C3").Value = 0.8787496612563173
$ws.Range("C5").Value = 2.431458940166964
$ws.Range("C6").Value = 1.447930496829564
$ws.Range("C8").Value = 0.6742451383204839
$ws.Range("C11").Value = 1.715791310593251
$ws.Range("C13").Value = 0.893498267486792
$ws.Range("C16").Value = 0.06579575777907465

# Column E (y_1_forecast)
$ws.Range("E3").Value = 1.013823151053028
$ws.Range("E4").Value = 4.356912452939454
$ws.Range("E5").Value = 5.259925231829876
$ws.Range("E7").Value = -0.563208905821222
$ws.Range("E8").Value = 1.713290556413583
$ws.Range("E11").Value = 1.687339605296501
$ws.Range("E12").Value = -2.079848588862143
$ws.Range("E13").Value = -1.194610791899986
$ws.Range("E14").Value = 0.8024032016000104
$ws.Range("E15").Value = 0.8023688159249032
$ws.Range("E16").Value = 0.2740865344839749
$ws.Range("E18").Value = -1.696610696428313
$ws.Range("E19").Value = -1.362365718491854
